# Update NATMI LR-pair sheet (Fgf17-Fgfr3) with new TPM-based computation.
# This adds a new "Resolving-Mac" sending/target cluster (in addition to the
# existing FAPs / MuSCs / ECs clusters) and refreshes all the derived
# expression-specificity metrics for every Sending-cluster x Target-cluster
# combination (rows 2-10), per columns:
#   A: Sending cluster            K: Receptor-expressing cells
#   B: Ligand symbol               L: Receptor detection rate
#   C: Receptor symbol             M: Receptor average expression value
#   D: Target cluster              N: Receptor total expression value
#   E: Ligand-expressing cells     O: Receptor derived specificity (avg)
#   F: Ligand detection rate       P: Receptor derived specificity (total)
#   G: Ligand average expression   Q: Edge average expression weight
#   H: Ligand total expression     R: Edge total expression weight
#   I: Ligand derived specificity (avg)   S: Edge average expression derived specificity
#   J: Ligand derived specificity (total) T: Edge total expression derived specificity

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    # A                B         C         D                 E  F                   G                    H                    I                     J                     K  L  M                    N          O                     P                     Q                   R                    S                     T
    @("FAPs",          "Fgf17", "Fgfr3", "ECs",           3, 1,                  1.430739666666667,  4.292218999999999,  0.7311024264480693, 0.7311024264480693, 3, 1, 5.751166666666666,  17.2535,  0.7405222614421495,  0.7405222614421495,  8.22842227961111,   74.05580051649999,  0.5413976221791671,   0.5413976221791671),
    @("FAPs",          "Fgf17", "Fgfr3", "FAPs",          3, 1,                  1.430739666666667,  4.292218999999999,  0.7311024264480693, 0.7311024264480693, 3, 1, 0.5698483333333333, 1.709545, 0.07337387367415998, 0.07337387367416,    0.8153046144838887, 7.337741530354998,  0.05364381708107248,  0.05364381708107249),
    @("FAPs",          "Fgf17", "Fgfr3", "Resolving-Mac", 3, 1,                  1.430739666666667,  4.292218999999999,  0.7311024264480693, 0.7311024264480693, 3, 1, 1.445350666666667,  4.336052, 0.1861038648836906,  0.1861038648836906,  2.067920531043111,  18.611284779388,    0.1360609871878298,   0.1360609871878298),
    @("Resolving-Mac", "Fgf17", "Fgfr3", "ECs",           2, 0.6666666666666666, 0.4223926666666666, 1.267178,            0.2158410161600821,  0.2158410161600822, 3, 1, 5.751166666666666,  17.2535,  0.7405222614421495,  0.7405222614421495,  2.429250624777778,  21.863255623,       0.1598350773988356,   0.1598350773988356),
    @("Resolving-Mac", "Fgf17", "Fgfr3", "FAPs",          2, 0.6666666666666666, 0.4223926666666666, 1.267178,            0.2158410161600821,  0.2158410161600822, 3, 1, 0.5698483333333333, 1.709545, 0.07337387367415998, 0.07337387367416,    0.2406997571122222, 2.16629781401,      0.01583709145343219,  0.0158370914534322),
    @("Resolving-Mac", "Fgf17", "Fgfr3", "Resolving-Mac", 2, 0.6666666666666666, 0.4223926666666666, 1.267178,            0.2158410161600821,  0.2158410161600822, 3, 1, 1.445350666666667,  4.336052, 0.1861038648836906,  0.1861038648836906,  0.6105055223617778, 5.494549701256,     0.04016884730781441,  0.04016884730781441),
    @("MuSCs",         "Fgf17", "Fgfr3", "ECs",           1, 0.3333333333333333, 0.1038296666666667, 0.311489,            0.05305655739184852, 0.05305655739184854,3, 1, 5.751166666666666,  17.2535,  0.7405222614421495,  0.7405222614421495,  0.5971417179444444, 5.3742754615,       0.03928956186414686,  0.03928956186414687),
    @("MuSCs",         "Fgf17", "Fgfr3", "FAPs",          1, 0.3333333333333333, 0.1038296666666667, 0.311489,            0.05305655739184852, 0.05305655739184854,3, 1, 0.5698483333333333, 1.709545, 0.07337387367415998, 0.07337387367416,    0.05916716250055555,0.5325044625049999, 0.003892965139655312, 0.003892965139655314),
    @("MuSCs",         "Fgf17", "Fgfr3", "Resolving-Mac", 1, 0.3333333333333333, 0.1038296666666667, 0.311489,            0.05305655739184852, 0.05305655739184854,3, 1, 1.445350666666667,  4.336052, 0.1861038648836906,  0.1861038648836906,  0.1500702779364445, 1.350632501428,     0.009874030388046353, 0.009874030388046356)
)

$startRow = 2
$r = $startRow
foreach ($rowValues in $rows) {
    $c = 1
    foreach ($val in $rowValues) {
        $ws.Cells.Item($r, $c).Value = $val
        $c++
    }
    $r++
}
